# refactor tension calculation and tests
# Applies to /tmp/work/before.xlsx:
#  - adds a new blank "Sheet1" worksheet after "Compression"
#  - reworks the "Compression" sheet: inserts new input columns
#    (Shear Modulus, inertia x/y, torsional/warping constants,
#    length/k for torsion) and adds a new "elastic torsional stress"
#    calculation block (columns AH:AN and AP)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Add the new trailing worksheet "Sheet1" (sheetId 4) after
#    "Compression", which is currently the last sheet.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

# ---------------------------------------------------------------
# 2) Rework the "Compression" sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Compression")
$ws.Activate()

# Insert the new columns (left to right so each subsequent insert
# position is expressed in the already-shifted column coordinates):
#   - 1 new column at B            (Shear Modulus Gpa)
#   - 4 new columns at D:G         (inertia x, inertia y, torsional
#                                    constant, warping constant)
#   - 2 new columns at O:P         (length torsion, k torsion)
$ws.Range("B:B").Insert()
$ws.Range("D:G").Insert()
$ws.Range("O:P").Insert()

# Populate the new header/value cells. The order in which brand new
# text values are assigned controls the order new shared-string
# entries are created, so the new torsional-stress header is set
# first, followed by the other new headers left-to-right.
$ws.Range("AH1").Value = "elastic torsional stress"
$ws.Range("D1").Value = "inertia x"
$ws.Range("E1").Value = "inertia y"
$ws.Range("F1").Value = "torsional constant"
$ws.Range("G1").Value = "warping constant"
$ws.Range("B1").Value = "Shear Modulus Gpa"
$ws.Range("O1").Value = "length torsion"
$ws.Range("P1").Value = "k torsion"

$ws.Range("B2").Value = 77
$ws.Range("D2").Value = 12.1
$ws.Range("E2").Value = 3.88
$ws.Range("F2").Value = 42
$ws.Range("G2").Value = 20.5
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

# Two of the pre-existing formulas were edited by hand to anchor a
# column reference absolutely.
$ws.Range("AA2").Formula = '=$C2/Z2'
$ws.Range("AE2").Formula = '=AD2*$J2'

# New "elastic torsional stress" block, columns AH:AN (+AP scratch
# column), reusing the existing header strings for the downstream
# steps shared with the other critical-stress blocks.
$ws.Range("AI1").Value = "Fy/Fe y"
$ws.Range("AJ1").Value = "critical stress y 1"
$ws.Range("AK1").Value = "critical stress y 2 "
$ws.Range("AL1").Value = "criticak stress"
$ws.Range("AM1").Value = "nominal strength"
$ws.Range("AN1").Value = "design strength"

$ws.Range("AH2").Formula = "=(PI()^2*A2*10^3*G2*10^9/(P2*O2*1000)^2+B2*10^3*F2*10^3)*1/(D2*10^6+E2*10^6)"
$ws.Range("AI2").Formula = '=$C2/AH2'
$ws.Range("AJ2").Formula = "=0.658^AI2*C2"
$ws.Range("AK2").Formula = "=0.877*AH2"
$ws.Range("AL2").Formula = "=IF(AI2<=2.25,AJ2,AK2)"
$ws.Range("AM2").Formula = '=AL2*$J2'
$ws.Range("AN2").Formula = "=AM2/1.67"
$ws.Range("AP2").Formula = "=(PI()^2*A2*10^3*G2*10^9/(P2*O2*1000)^2+B2*10^3*F2*10^3)"

# Best-effort cosmetic column widths for the newly inserted columns.
$ws.Range("B:B").ColumnWidth = 18.17
$ws.Range("D:G").ColumnWidth = 18.1
$ws.Range("AH:AH").ColumnWidth = 39.04
$ws.Range("AI:AI").ColumnWidth = 13.23
$ws.Range("AJ:AJ").ColumnWidth = 17.23

# Final selection/view state recorded in the saved file.
$ws.Range("AP2").Select()
